$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G10").Value = "Marina.youhana@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"
$ws.Range("G15").Value = "nancy.abdelshafy@med.asu.edu.eg, hana.amr@med.asu.edu.eg"
$ws.Range("G18").Value = "Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G19").Value = "Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G20").Value = "Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G24").Value = "servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G32").Value = "Marina.youhana@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"
$ws.Range("G37").Value = "nancy.abdelshafy@med.asu.edu.eg, hana.amr@med.asu.edu.eg"
$ws.Range("G40").Value = "Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G41").Value = "Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G42").Value = "Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G46").Value = "hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg"
$ws.Range("G52").Value = "Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G54").Value = "Marina.youhana@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G62").Value = "yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G63").Value = "Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("G68").Value = "hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg"
$ws.Range("G74").Value = "Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G76").Value = "Marina.youhana@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G81").Value = "user@user.com, enas.omran@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg"
$ws.Range("G84").Value = "yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G85").Value = "Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("G90").Value = "shaimaa.ahmed@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G96").Value = "aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G98").Value = "Marina.youhana@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G106").Value = "Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G107").Value = "Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G108").Value = "Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G112").Value = "shaimaa.ahmed@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G118").Value = "aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G120").Value = "Marina.youhana@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G128").Value = "Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G129").Value = "Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G130").Value = "Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G134").Value = "Amira.Sobhy@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G142").Value = "basma.hamed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg"
$ws.Range("G150").Value = "Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G156").Value = "alshimaa.atef@med.asu.edu.egm, majorelle.magdy@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G164").Value = "basma.hamed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg"
$ws.Range("G172").Value = "Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
